# Apply updated Fitness values (column C) for run_21 log sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-11 (Generation 0-9): Fitness 7310 -> 7318
$ws.Range("C2:C11").Value = 7318

# Rows 44-106 (Generation 42-104): Fitness 7295/7293 -> 7310
$ws.Range("C44:C106").Value = 7310
